$wb = $excel.ActiveWorkbook

# Helper: write a value into a cell, forcing it to be stored as TEXT
# (shared string) even when the text looks like a plain number, mirroring
# how the source data (numeric-looking labels) is stored in the workbook.
# A leading apostrophe forces Excel to treat the entry as text; resetting
# the Style back to "Normal" afterwards strips the quote-prefix number
# format so the cell keeps the workbook's default (un-styled) cell format.
function Set-TextValue {
    param($Range, [string]$Text)
    $Range.Value = "'" + $Text
    $Range.Style = "Normal"
}

# --- Sheet 2: Restricciones_del_lider ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(2, 1).Value = "0.049999999999998934 - x + y"
Set-TextValue $ws2.Cells.Item(2, 2) "-0.049999999999998934"
Set-TextValue $ws2.Cells.Item(2, 4) "0.4"

# --- Sheet 3: Restricciones_del_follower ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(2, 1).Value = "-28.85 + x + y"
Set-TextValue $ws3.Cells.Item(2, 2) "8.85"
Set-TextValue $ws3.Cells.Item(2, 4) "0.55"
Set-TextValue $ws3.Cells.Item(2, 6) "6.5"

$ws3.Cells.Item(3, 1).Value = "14.4 - y"
Set-TextValue $ws3.Cells.Item(3, 2) "-14.4"
Set-TextValue $ws3.Cells.Item(3, 4) "0.45"
Set-TextValue $ws3.Cells.Item(3, 6) "1.1"

$ws3.Cells.Item(4, 1).Value = "-25.6 + y"
Set-TextValue $ws3.Cells.Item(4, 2) "-5.6"
Set-TextValue $ws3.Cells.Item(4, 4) "0.55"
Set-TextValue $ws3.Cells.Item(4, 5) "0"
Set-TextValue $ws3.Cells.Item(4, 6) "1.6"

# --- Sheet 4: Punto_modificado ---
$ws4 = $wb.Worksheets.Item(4)
Set-TextValue $ws4.Cells.Item(2, 1) "14.45"
Set-TextValue $ws4.Cells.Item(2, 2) "14.4"

# --- Sheet 5: Vector_bf ---
$ws5 = $wb.Worksheets.Item(5)
Set-TextValue $ws5.Cells.Item(2, 1) "-53.65"

# --- Sheet 6: Vector_BF ---
$ws6 = $wb.Worksheets.Item(6)
Set-TextValue $ws6.Cells.Item(2, 1) "-28.5"
Set-TextValue $ws6.Cells.Item(3, 1) "-9.200000000000001"
